# Sync file from Google Drive
#
# 1. Drop the "Right_Right_BusStopCode" column (S) from each NextBus sheet;
#    the former "Bus Stop Description" column (T) shifts left into S.
# 2. Refresh the EstimatedTimeOfArrival (F) timestamps with the newer pull.
# 3. One row's TypeOfBus flips from SD -> DD (NextBus2!L9).

$wb = $excel.ActiveWorkbook

# --- 1. Remove column S (Right_Right_BusStopCode) on every sheet ---------
foreach ($sheetName in @("NextBus1", "NextBus2", "NextBus3")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("S1").EntireColumn.Delete()
}

# --- 2. Updated EstimatedTimeOfArrival values -----------------------------
$timeUpdates = @{
    "NextBus1" = @{
        "F2"  = 45684.97600694445
        "F3"  = 45684.9807175926
        "F4"  = 45684.97912037037
        "F5"  = 45684.97212962963
        "F6"  = 45684.97432870371
        "F7"  = 45684.98070601852
        "F8"  = 45684.97855324074
        "F9"  = 45684.97578703704
        "F10" = 45684.9794675926
        "F11" = 45684.97868055556
        "F13" = 45684.97733796296
        "F14" = 45684.97723379629
    }
    "NextBus2" = @{
        "F2"  = 45684.98638888889
        "F3"  = 45684.9887962963
        "F4"  = 45684.98086805556
        "F5"  = 45684.98016203703
        "F7"  = 45684.98082175926
        "F8"  = 45684.98466435185
        "F9"  = 45684.98666666666
        "F11" = 45684.98451388889
        "F12" = 45684.98135416667
        "F13" = 45684.98634259259
    }
    "NextBus3" = @{
        "F3" = 45684.99162037037
        "F4" = 45684.98528935185
        "F5" = 45684.98553240741
        "F6" = 45684.99459490741
        "F7" = 45684.99128472222
        "F8" = 45684.99364583333
        "F9" = 45684.99446759259
    }
}

foreach ($sheetName in $timeUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $timeUpdates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}

# --- 3. NextBus2!L9 TypeOfBus: SD -> DD -----------------------------------
$ws2 = $wb.Worksheets.Item("NextBus2")
$ws2.Range("L9").Value = "DD"
